$d = $word.ActiveDocument

# 1. Apply the "DejaVu Sans" font (ascii/hAnsi via Name, cs via NameBi) to
#    every paragraph in the document - this stamps w:rFonts onto every run's
#    rPr as well as onto the paragraph-mark rPr stored in each pPr.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $r.Font.Name = "DejaVu Sans"
    $r.Font.NameBi = "DejaVu Sans"
}

# 2. Update the cached MERGEFIELD result for "City" from the old value to
#    the new accommodation city, everywhere it appears in the document.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("416071059", $false, $false, $false, $false, $false, $true, 1, $false, "VICO DEL GARGANO", 2)
